# Apply the update to the "Assets" sheet in Config.xlsx:
# add three new Generic Asset rows (Name/Value columns mirror each other,
# Description column references the "Generic Asset" shared string) and
# move the active selection to B7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Assets")

$ws.Range("A3").Value = "FolderLocation_LogFiles"
$ws.Range("B3").Value = "FolderLocation_LogFiles"
$ws.Range("C3").Value = "Generic Asset"

$ws.Range("A4").Value = "FolderLocation_Screenshots"
$ws.Range("B4").Value = "FolderLocation_Screenshots"
$ws.Range("C4").Value = "Generic Asset"

$ws.Range("A5").Value = "FolderLocation_UserInterface"
$ws.Range("B5").Value = "FolderLocation_UserInterface"
$ws.Range("C5").Value = "Generic Asset"

$ws.Range("B7").Select() | Out-Null
